# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 265 (pushes existing rows 265-272 down to 266-273)
# and populate it with the new "Elegant Lady" Durazno record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 265, shifting rows 265:272 down to 266:273.
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row 265 with the new record.
$ws.Cells.Item(265, 1).Value = 10
$ws.Cells.Item(265, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(265, 3).Value = "La Araucanía"
$ws.Cells.Item(265, 4).Value = 44628
$ws.Cells.Item(265, 5).Value = 9
$ws.Cells.Item(265, 6).Value = "Fruta"
$ws.Cells.Item(265, 7).Value = 100103
$ws.Cells.Item(265, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(265, 9).Value = 100103004
$ws.Cells.Item(265, 10).Value = "Durazno"
$ws.Cells.Item(265, 11).Value = "Elegant Lady"
$ws.Cells.Item(265, 12).Value = "Primera"
$ws.Cells.Item(265, 13).Value = 100
$ws.Cells.Item(265, 14).Value = 17000
$ws.Cells.Item(265, 15).Value = 17000
$ws.Cells.Item(265, 16).Value = 17000
$ws.Cells.Item(265, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(265, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(265, 19).Value = 944
$ws.Cells.Item(265, 20).Value = 18
